$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated crypto price/volume figures (and two swapped row-pairs) per the
# "Wed Dec 27 08:27:16 UTC 2023" GitHub Actions data refresh.
# Each entry is Row -> @{ Col = NewValue } for every cell that changed.
$updates = @{
    2 = @{ 'D'='42.509.57'; 'E'='  -0.98%  ' }
    3 = @{ 'D'='2.234.96'; 'E'='  -0.10%  ' }
    4 = @{ 'E'='  +0.11%  ' }
    5 = @{ 'D'='111.93'; 'E'='  -0.46%  ' }
    6 = @{ 'D'='294.10'; 'E'='  +5.51%  ' }
    7 = @{ 'D'='0.621'; 'E'='  -0.95%  ' }
    8 = @{ 'E'='  -0.41%  ' }
    9 = @{ 'D'='0.600'; 'E'='  -1.03%  ' }
    10 = @{ 'E'='  -5.64%  ' }
    11 = @{ 'D'='0.0919'; 'E'='  -0.84%  ' }
    12 = @{ 'D'='54.10'; 'E'='  +0.12%  ' }
    13 = @{ 'E'='  -2.25%  ' }
    14 = @{ 'E'='  +22.22%  ' }
    15 = @{ 'E'='  -1.31%  ' }
    16 = @{ 'D'='14.97'; 'E'='  -2.02%  ' }
    17 = @{ 'D'='2.569.34'; 'E'='  -0.23%  ' }
    18 = @{ 'D'='2.238.28'; 'E'='  -0.05%  ' }
    19 = @{ 'D'='42.483.63'; 'E'='  -1.03%  ' }
    20 = @{ 'B'='ShibaInu'; 'C'='https://coinranking.com/coin/xz24e0BjL+shibainu-shib'; 'D'='0.0000106'; 'E'='  -1.13%  ' }
    21 = @{ 'B'='Uniswap'; 'C'='https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'; 'D'='7.18'; 'E'='  +6.56%  ' }
    22 = @{ 'D'='74.04'; 'E'='  +2.91%  ' }
    23 = @{ 'E'='  +10.29%  ' }
    24 = @{ 'D'='2.37'; 'E'='  +1.51%  ' }
    25 = @{ 'D'='244.43'; 'E'='  +5.69%  ' }
    26 = @{ 'D'='8.90'; 'E'='  -3.44%  ' }
    27 = @{ 'D'='1.00'; 'E'='  -0.94%  ' }
    28 = @{ 'D'='11.40'; 'E'='  -5.97%  ' }
    29 = @{ 'E'='  -1.77%  ' }
    30 = @{ 'D'='175.24'; 'E'='  +0.92%  ' }
    31 = @{ 'D'='37.10'; 'E'='  -7.92%  ' }
    32 = @{ 'B'='EthereumClassic'; 'C'='https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'; 'D'='21.61'; 'E'='  +2.51%  ' }
    33 = @{ 'B'='WEMIXToken'; 'C'='https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'; 'D'='3.12'; 'E'='  -4.31%  ' }
    34 = @{ 'D'='0.0879'; 'E'='  -2.45%  ' }
    35 = @{ 'D'='5.69'; 'E'='  +2.42%  ' }
    36 = @{ 'D'='5.00'; 'E'='  +7.76%  ' }
    37 = @{ 'B'='NEARProtocol'; 'C'='https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'; 'D'='4.22'; 'E'='  -1.54%  ' }
    38 = @{ 'B'='Stellar'; 'C'='https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'; 'D'='0.127'; 'E'='  -0.31%  ' }
    39 = @{ 'D'='0.0372'; 'E'='  +0.33%  ' }
    40 = @{ 'D'='0.104'; 'E'='  -1.88%  ' }
    41 = @{ 'D'='2.40'; 'E'='  -6.60%  ' }
    42 = @{ 'D'='70.96'; 'E'='  +0.17%  ' }
    43 = @{ 'D'='0.228'; 'E'='  -1.10%  ' }
    44 = @{ 'E'='  -0.18%  ' }
    45 = @{ 'D'='12.33'; 'E'='  -6.01%  ' }
    46 = @{ 'E'='  -2.18%  ' }
    47 = @{ 'D'='5.44'; 'E'='  -3.61%  ' }
    48 = @{ 'D'='1.29'; 'E'='  +0.58%  ' }
    49 = @{ 'B'='FraxShare'; 'C'='https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'; 'D'='8.51'; 'E'='  +1.49%  ' }
    50 = @{ 'B'='Aave'; 'C'='https://coinranking.com/coin/ixgUfzmLR+aave-aave'; 'D'='102.18'; 'E'='  +2.04%  ' }
    51 = @{ 'D'='0.0976'; 'E'='  -0.94%  ' }
}

foreach ($row in $updates.Keys) {
    foreach ($col in $updates[$row].Keys) {
        $value = $updates[$row][$col]
        $cell = $ws.Range("$col$row")
        if ($col -eq "D" -or $col -eq "E") {
            # Force text so Excel does not coerce numeric-looking
            # strings (e.g. "111.93", "1.00") into floating point
            # numbers, and does not trim the padded "  +x.xx%  " cells.
            $cell.NumberFormat = "@"
            $cell.Value = $value
            $cell.ClearFormats()
        } else {
            $cell.Value = $value
        }
    }
}
